$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

# "Recorded By" (column G) values were reordered from "System, <email>"
# to "<email>, System" throughout the sheet.
$oldValue = "System, dnasr281@gmail.com"
$newValue = "dnasr281@gmail.com, System"

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1
$col = 7  # column G

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, $col)
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value2 = $newValue
    }
}
